# Updates the DPM integration testfixture with new hierarchy node labels.
#
# Changes applied (per the commit's canonical-OOXML diff):
#   - CodeSchemes!A2   : new scheme ID (UUID)
#   - Codes!A2, Codes!A3 : new code IDs (UUID)
#   - Extensions!A2, Extensions!A3 : new extension IDs (UUID)
#   - Members_HIER!A2  : new member ID (UUID)
#   - Members_HIER!D2  : "Member (fi, label)" -> "Hierarchy node (fi, label)"
#   - Column A autofit widths shifted slightly on every sheet (and column D
#     on Members_HIER) as a side effect of the regenerated fixture data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: CodeSchemes
# ---------------------------------------------------------------------
$wsCodeSchemes = $wb.Worksheets.Item("CodeSchemes")
$wsCodeSchemes.Range("A2").Value = "a68e6334-aedf-4a28-9964-057446f62d8c"
$wsCodeSchemes.Columns.Item(1).ColumnWidth = 31.142857142857142

# ---------------------------------------------------------------------
# Sheet 2: Codes
# ---------------------------------------------------------------------
$wsCodes = $wb.Worksheets.Item("Codes")
$wsCodes.Range("A2").Value = "b8de6ed6-4461-433f-bd28-d117499ebf41"
$wsCodes.Range("A3").Value = "b8887fac-bca8-408f-8f84-6755e0e4370e"
$wsCodes.Columns.Item(1).ColumnWidth = 35.57142857142857

# ---------------------------------------------------------------------
# Sheet 3: Extensions
# ---------------------------------------------------------------------
$wsExtensions = $wb.Worksheets.Item("Extensions")
$wsExtensions.Range("A2").Value = "4ffe9ff2-c8d2-4595-aa32-1155ac4598eb"
$wsExtensions.Range("A3").Value = "bc5cc9ce-7f49-4727-a2a6-ac7c390398b0"
$wsExtensions.Columns.Item(1).ColumnWidth = 27.857142857142858

# ---------------------------------------------------------------------
# Sheet 5: Members_HIER
# ---------------------------------------------------------------------
$wsMembersHier = $wb.Worksheets.Item("Members_HIER")
$wsMembersHier.Range("A2").Value = "29ce4f17-4468-436b-b4d8-300235d45e69"
$wsMembersHier.Range("D2").Value = "Hierarchy node (fi, label)"
$wsMembersHier.Columns.Item(1).ColumnWidth = 35.57142857142857
$wsMembersHier.Columns.Item(4).ColumnWidth = 18.0
